$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Average_Class_Volatility" column (K) was dropped from the analysis.
# Deleting it shifts Beta / Correlation_with_Market / R_Squared / Data_Points
# one column to the left (L->K, M->L, N->M, O->N) and updates the used range
# from A1:O10 down to A1:N10.
$ws.Columns.Item(11).Delete()

# Row 2: AAPL_STOCK
$ws.Range("B2").Value = 0.0006
$ws.Range("C2").Value = 0.0175
$ws.Range("D2").Value = 0.1584
$ws.Range("E2").Value = 0.278
$ws.Range("F2").Value = 0.5697
$ws.Range("G2").Value = 0.9491
$ws.Range("H2").Value = 14.0461
$ws.Range("I2").Value = -0.0925
$ws.Range("J2").Value = 0.1533
$ws.Range("K2").Value = 1.184198104117674
$ws.Range("L2").Value = 0.7028000837994652
$ws.Range("M2").Value = 0.4939279577885354
$ws.Range("N2").Value = 501

# Row 3: GOOGL_STOCK
$ws.Range("B3").Value = 0.0011
$ws.Range("C3").Value = 0.0186
$ws.Range("D3").Value = 0.2693
$ws.Range("E3").Value = 0.295
$ws.Range("F3").Value = 0.913
$ws.Range("G3").Value = -0.1862
$ws.Range("H3").Value = 5.3597
$ws.Range("I3").Value = -0.0951
$ws.Range("J3").Value = 0.1022
$ws.Range("K3").Value = 1.09040109559776
$ws.Range("L3").Value = 0.6099254016694213
$ws.Range("M3").Value = 0.3720089956016049
$ws.Range("N3").Value = 501

# Row 4: TSLA_STOCK
$ws.Range("B4").Value = 0.0015
$ws.Range("C4").Value = 0.0401
$ws.Range("D4").Value = 0.381
$ws.Range("E4").Value = 0.6363
$ws.Range("F4").Value = 0.5988
$ws.Range("G4").Value = 0.6843
$ws.Range("H4").Value = 4.6366
$ws.Range("I4").Value = -0.1543
$ws.Range("J4").Value = 0.2269
$ws.Range("K4").Value = 2.260707725581899
$ws.Range("L4").Value = 0.5862692072607113
$ws.Range("M4").Value = 0.3437115833821028
$ws.Range("N4").Value = 501

# Row 5: SPY_STOCK
$ws.Range("B5").Value = 0.0009
$ws.Range("C5").Value = 0.0104
$ws.Range("D5").Value = 0.216
$ws.Range("E5").Value = 0.165
$ws.Range("F5").Value = 1.3088
$ws.Range("G5").Value = 1.1333
$ws.Range("H5").Value = 23.1748
$ws.Range("I5").Value = -0.0585
$ws.Range("J5").Value = 0.105
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 501

# Row 6: QQQ_STOCK
$ws.Range("B6").Value = 0.001
$ws.Range("C6").Value = 0.0132
$ws.Range("D6").Value = 0.2485
$ws.Range("E6").Value = 0.2098
$ws.Range("F6").Value = 1.1843
$ws.Range("G6").Value = 0.8108
$ws.Range("H6").Value = 14.3248
$ws.Range("I6").Value = -0.0621
$ws.Range("J6").Value = 0.12
$ws.Range("K6").Value = 1.219462404147546
$ws.Range("L6").Value = 0.9588425327119318
$ws.Range("M6").Value = 0.9193790025374319
$ws.Range("N6").Value = 501

# Row 7: GLD_COMMODITY
$ws.Range("B7").Value = 0.0012
$ws.Range("C7").Value = 0.01
$ws.Range("D7").Value = 0.2931
$ws.Range("E7").Value = 0.1591
$ws.Range("F7").Value = 1.8417
$ws.Range("G7").Value = -0.043
$ws.Range("H7").Value = 1.0061
$ws.Range("I7").Value = -0.0357
$ws.Range("J7").Value = 0.037
$ws.Range("K7").Value = 0.09573965558166579
$ws.Range("L7").Value = 0.09926356379498978
$ws.Range("M7").Value = 0.009853255097282005
$ws.Range("N7").Value = 501

# Row 8: SLV_COMMODITY
$ws.Range("B8").Value = 0.0011
$ws.Range("C8").Value = 0.0176
$ws.Range("D8").Value = 0.2707
$ws.Range("E8").Value = 0.2789
$ws.Range("F8").Value = 0.9707
$ws.Range("G8").Value = 0.0029
$ws.Range("H8").Value = 1.3797
$ws.Range("I8").Value = -0.0652
$ws.Range("J8").Value = 0.0639
$ws.Range("K8").Value = 0.4633620344859677
$ws.Range("L8").Value = 0.2741491321802776
$ws.Range("M8").Value = 0.07515774667519934
$ws.Range("N8").Value = 501

# Row 9: USO_COMMODITY
$ws.Range("B9").Value = 0.0003
$ws.Range("C9").Value = 0.0187
$ws.Range("D9").Value = 0.0668
$ws.Range("E9").Value = 0.2973
$ws.Range("F9").Value = 0.2248
$ws.Range("G9").Value = -0.2112
$ws.Range("H9").Value = 1.4209
$ws.Range("I9").Value = -0.0807
$ws.Range("J9").Value = 0.0689
$ws.Range("K9").Value = 0.3273858226395131
$ws.Range("L9").Value = 0.1817129300185401
$ws.Range("M9").Value = 0.03301958893592283
$ws.Range("N9").Value = 501

# Row 10: DBA_COMMODITY
$ws.Range("B10").Value = 0.0007
$ws.Range("C10").Value = 0.0096
$ws.Range("D10").Value = 0.1826
$ws.Range("E10").Value = 0.1529
$ws.Range("F10").Value = 1.1946
$ws.Range("G10").Value = -0.3783
$ws.Range("H10").Value = 1.8715
$ws.Range("I10").Value = -0.0486
$ws.Range("J10").Value = 0.0303
$ws.Range("K10").Value = 0.1779771087446846
$ws.Range("L10").Value = 0.1920722136958383
$ws.Range("M10").Value = 0.03689173527401976
$ws.Range("N10").Value = 501

